$d = $word.ActiveDocument

# List of (old, new) text replacements, applied in document order so that
# a newly-inserted value never collides with a not-yet-processed search term.
$replacements = @(
    @("14÷5=2, 4",  "10÷3=3, 1"),
    @("93÷5=18, 3", "63÷8=7, 7"),
    @("42÷3=14, 0", "34÷2=17, 0"),
    @("96÷8=12, 0", "46÷3=15, 1"),
    @("12÷7=1, 5",  "78÷9=8, 6"),
    @("61÷6=10, 1", "56÷5=11, 1"),
    @("56÷4=14, 0", "56÷7=8, 0"),
    @("60÷8=7, 4",  "32÷6=5, 2"),
    @("20÷5=4, 0",  "19÷4=4, 3"),
    @("42÷7=6, 0",  "25÷6=4, 1"),
    @("57÷8=7, 1",  "46÷7=6, 4"),
    @("70÷4=17, 2", "33÷6=5, 3"),
    @("56÷6=9, 2",  "73÷4=18, 1"),
    @("96÷5=19, 1", "86÷3=28, 2"),
    @("67÷3=22, 1", "94÷2=47, 0"),
    @("58÷7=8, 2",  "67÷3=22, 1"),
    @("50÷7=7, 1",  "61÷8=7, 5"),
    @("29÷5=5, 4",  "43÷2=21, 1"),
    @("47÷9=5, 2",  "42÷6=7, 0"),
    @("46÷6=7, 4",  "49÷8=6, 1"),
    @("44÷5=8, 4",  "17÷6=2, 5"),
    @("60÷4=15, 0", "21÷6=3, 3"),
    @("34÷3=11, 1", "41÷6=6, 5"),
    @("49÷4=12, 1", "36÷6=6, 0"),
    @("75÷6=12, 3", "80÷9=8, 8")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
